# "+compared data of all models"
# Update the freshly-computed comparison metrics on the "in" sheet (row 4 =
# per-model error stats, row 5 = raw counts) and refresh the view state that
# Excel records after the author finished reviewing the numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 4: corrected metric values (were stale/out-of-scale values from a
# previous run before all models' outputs were compared on equal footing) ---
$ws.Range("X4").Value = 68.4252
$ws.Range("Y4").Value = 60.21412
$ws.Range("Z4").Value = 58.423
$ws.Range("AD4").Value = 57.899
$ws.Range("AI4").Value = 55.787617
$ws.Range("AJ4").Value = 54.784353
$ws.Range("AN4").Value = 72.41212
$ws.Range("AO4").Value = 69.2132

# --- Row 5: corrected count ---
$ws.Range("AN5").Value = 752

# --- View state: switch the printed page to portrait and leave the window
# zoomed/scrolled to where the author ended up after comparing the models ---
$ws.PageSetup.Orientation = 1
$excel.ActiveWindow.Zoom = 86
$ws.Range("AN5").Select() | Out-Null
